$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update quantity/amount text away from "case"/"lbs" wording
$ws.Range("B1").Value = "1 p1fett 4E"
$ws.Range("B2").Value = "2 p2papp 5E"

# Move the active selection to B2
$ws.Range("B2").Select()
